$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'33.947.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.771.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.19%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.36%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'224.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.35%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.546"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.68%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.34%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'31.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.90%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0928"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.10%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.024.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.21%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +7.02%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.779.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.79%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'33.943.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.03%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.619"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.86%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'68.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.38%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'250.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.37%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -2.31%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.55%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'155.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.83%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'16.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.41%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'6.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.98%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.76%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.34%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.0509"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.49%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.00%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.16%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.62%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.443.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -6.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.24%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.620"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.97%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.00%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'82.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.41%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.89%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.879"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.07%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -4.95%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.43%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.66%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.924.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.02%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.42%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'FraxShare"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'5.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.07%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'11.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.94%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Quant"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'96.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'BitcoinSV"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'49.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -6.09%  "
$ws.Range("E51").Style = "Normal"
